# Update the "cryptos" price/volume snapshot with refreshed data
# (GitHub Actions scheduled scrape), matching commit:
# "Updated cryptos list on Fri Jul 19 11:12:13 UTC 2024 with GitHub Actions"
#
# Column D = Price (text), Column E = Volume(1h) change (text, % with padding).
# Values that look like plain decimal numbers (single "." , e.g. "162.00") are
# written with a leading apostrophe so Excel keeps storing them as text,
# exactly like the original workbook (prices such as "64.008.94" already
# contain two dots and are never auto-converted, so they are left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.008.94'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '3.401.54'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''571.57'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = '''162.00'
$ws.Range("E6").Value = '  +1.98%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.400.65'
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").Value = '''0.548'
$ws.Range("E9").Value = '  -5.01%  '
$ws.Range("D10").Value = '''7.31'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("E12").Value = '  -4.14%  '
$ws.Range("E13").Value = '  -1.12%  '
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("E15").Value = '  -2.35%  '
$ws.Range("E16").Value = '  -1.33%  '
$ws.Range("D17").Value = '64.050.19'
$ws.Range("E17").Value = '  -0.84%  '
$ws.Range("D18").Value = '3.430.06'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").Value = '''13.42'
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").Value = '''373.25'
$ws.Range("E21").Value = '  -1.35%  '
$ws.Range("D22").Value = '''7.78'
$ws.Range("E22").Value = '  -1.67%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '''70.32'
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("E25").Value = '  -3.54%  '
$ws.Range("E26").Value = '  -4.27%  '
$ws.Range("D27").Value = '''9.46'
$ws.Range("E27").Value = '  -4.58%  '
$ws.Range("D28").Value = '''0.178'
$ws.Range("E28").Value = '  -0.70%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = '''6.06'
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("E32").Value = '  -0.79%  '
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").Value = '''22.75'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("D35").Value = '''7.02'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("E36").Value = '  -6.47%  '
$ws.Range("D37").Value = '''159.25'
$ws.Range("E37").Value = '  -1.12%  '
$ws.Range("D38").Value = '''0.851'
$ws.Range("E38").Value = '  +7.22%  '
$ws.Range("D39").Value = '''1.82'
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("D40").Value = '''0.0723'
$ws.Range("E40").Value = '  -2.65%  '
$ws.Range("D41").Value = '''25.75'
$ws.Range("E41").Value = '  -1.75%  '
$ws.Range("D42").Value = '''42.77'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '''6.43'
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = '2.722.75'
$ws.Range("E44").Value = '  -5.50%  '
$ws.Range("D45").Value = '''25.68'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("E46").Value = '  -3.71%  '
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").Value = '''2.40'
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("D49").Value = '''329.03'
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("E50").Value = '  -2.97%  '
$ws.Range("E51").Value = '  -2.06%  '
